$d = $word.ActiveDocument

# 1. Append the extra sentence to the existing run's text.
$d.Content.Find.Execute(
    "今天学习了分支管理 创建了一个dev分支",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "今天学习了分支管理 创建了一个dev分支 使用git创建分支简单又快速", 2)

# 2. Split the paragraph right before the "_GoBack" bookmark so the
#    bookmark ends up alone in a new (empty) trailing paragraph.
$bm = $d.Bookmarks("_GoBack")
$r = $d.Range($bm.Start, $bm.Start)
$r.Text = [char]13
